$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the two "parallelogram" question/answer rows with two new physics
# formula questions (text content change only; styles/number format stay).
$ws.Range("A8").Value = "Назовите формулу силы тяжести"
$ws.Range("B8").Value = "F = mg"
$ws.Range("A9").Value = "Назовите формулу закона Ома"
$ws.Range("B9").Value = "I = U/R"

# Move the active selection to B9.
$ws.Range("B9").Select()

# Add a small empty text box near A12/B12 (as in the source file).
$shp = $ws.Shapes.AddTextbox(1, 675.75, 261.75, 16.5, 13.561181102362205)
$shp.Name = "TextBox 2"
$shp.TextFrame.Characters.Font.Size = 11
